# Add new "I0" and "IF" columns (I and J) to the worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: new columns I1/IF1 (I1 = "I0", J1 = "IF"), styled like the
# existing header cell H1 ("IP").
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Per-row data for columns I (I0) and J (IF), rows 2..78.
$data = @(
    @(2,8,8),
    @(3,2,2),
    @(4,7,7),
    @(5,7,7),
    @(6,8,9),
    @(7,10,10),
    @(8,5,5),
    @(9,8,8),
    @(10,9,9),
    @(11,7,7),
    @(12,8,8),
    @(13,7,7),
    @(14,5,6),
    @(15,8,8),
    @(16,8,8),
    @(17,7,7),
    @(18,7,7),
    @(19,9,9),
    @(20,7,7),
    @(21,8,8),
    @(22,7,7),
    @(23,7,7),
    @(24,6,7),
    @(25,8,8),
    @(26,7,7),
    @(27,6,6),
    @(28,7,7),
    @(29,6,6),
    @(30,9,9),
    @(31,7,7),
    @(32,7,7),
    @(33,8,8),
    @(34,9,9),
    @(35,9,9),
    @(36,8,8),
    @(37,7,7),
    @(38,7,7),
    @(39,8,8),
    @(40,9,9),
    @(41,7,7),
    @(42,7,7),
    @(43,7,7),
    @(44,7,7),
    @(45,7,7),
    @(46,7,7),
    @(47,7,7),
    @(48,7,7),
    @(49,7,7),
    @(50,8,8),
    @(51,7,7),
    @(52,9,9),
    @(53,7,7),
    @(54,7,8),
    @(55,8,8),
    @(56,7,9),
    @(57,7,7),
    @(58,8,8),
    @(59,7,8),
    @(60,9,9),
    @(61,8,8),
    @(62,7,8),
    @(63,7,7),
    @(64,9,9),
    @(65,7,7),
    @(66,7,7),
    @(67,7,8),
    @(68,7,8),
    @(69,7,7),
    @(70,9,9),
    @(71,8,8),
    @(72,7,7),
    @(73,7,7),
    @(74,7,8),
    @(75,10,10),
    @(76,3,3),
    @(77,5,5),
    @(78,4,4)
)

foreach ($entry in $data) {
    $row = $entry[0]
    $iVal = $entry[1]
    $jVal = $entry[2]
    $ws.Cells.Item($row, 9).Value = $iVal
    $ws.Cells.Item($row, 10).Value = $jVal
}
